$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.011.65'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  -2.95%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.924.30'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -3.62%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.42'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -1.08%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.15'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -4.84%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.506'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -1.98%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.923.02'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  -3.55%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.93'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +3.24%  '
$ws.Range('E11').Value = '  -4.15%  '
$ws.Range('E12').Value = '  -3.38%  '
$ws.Range('E13').Value = '  -3.22%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.70'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -5.38%  '
$ws.Range('E15').Value = '  +0.29%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.406.95'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -3.66%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.951.40'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -2.97%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.79'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  -4.02%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.922.03'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -3.70%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.98'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -4.89%  '
$ws.Range('E21').Value = '  -4.41%  '
$ws.Range('E22').Value = '  -1.97%  '
$ws.Range('E23').Value = '  -4.63%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.61'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -3.04%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.87'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('E26').Value = '  -2.37%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.96'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  +0.00%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.24'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  -3.31%  '
$ws.Range('E32').Value = '  -3.33%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.62'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -3.39%  '
$ws.Range('E34').Value = '  -3.12%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0877'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  +1.35%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.02'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('E37').Value = '  -4.55%  '
$ws.Range('E38').Value = '  -5.26%  '
$ws.Range('E39').Value = '  -0.99%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.51'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -1.73%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.02'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  -4.50%  '
$ws.Range('E42').Value = '  -4.36%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.299'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -1.64%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.64'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -1.36%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '379.65'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -3.57%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0350'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  -3.00%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.701.49'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.79%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.52'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +0.23%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.02'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +2.50%  '
$ws.Range('E51').Value = '  -1.68%  '
